$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp footer (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 12:46"

# Rows 12 & 13 swap: A Coruna now sorts above Alacant/Alicante, with refreshed counts;
# Alacant/Alicante keeps its previous totals.
$ws.Range("A12").Value = "A Coruña"
$ws.Range("B12").Value = 634
$ws.Range("C12").Value = 19
$ws.Range("D12").Value = 621
$ws.Range("E12").Value = 13

$ws.Range("A13").Value = "Alacant/Alicante"
$ws.Range("B13").Value = 628
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 555
$ws.Range("E13").Value = 61

# Pontevedra (row 17): refreshed totals/recovered counts
$ws.Range("B17").Value = 494
$ws.Range("D17").Value = 491

# Ourense (row 41): refreshed totals/recovered counts
$ws.Range("B41").Value = 152
$ws.Range("D41").Value = 150

# Rows 45-48: Lugo now sorts above Zamora (between Cuenca and Zamora), with refreshed
# counts; Zamora, Almeria and Palencia shift down a row but keep their previous totals.
$ws.Range("A45").Value = "Lugo"
$ws.Range("B45").Value = 97
$ws.Range("C45").Value = 19
$ws.Range("D45").Value = 94
$ws.Range("E45").Value = 3

$ws.Range("A46").Value = "Zamora"
$ws.Range("B46").Value = 90
$ws.Range("C46").Value = 5
$ws.Range("D46").Value = 81
$ws.Range("E46").Value = 4

$ws.Range("A47").Value = "Almeria"
$ws.Range("B47").Value = 74
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 72
$ws.Range("E47").Value = 2

$ws.Range("A48").Value = "Palencia"
$ws.Range("B48").Value = 64
$ws.Range("C48").Value = 2
$ws.Range("D48").Value = 62
$ws.Range("E48").Value = 0
